$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5782.357
$ws.Range("J43").Value = 6268.364
$ws.Range("L43").Value = 6268.364
$ws.Range("N43").Value = -6406.364

$ws.Range("H45").Value = 4189.75
$ws.Range("I45").Value = 580
$ws.Range("J45").Value = 15019
$ws.Range("K45").Value = 1740
$ws.Range("L45").Value = 45057
$ws.Range("M45").Value = -1548
$ws.Range("N45").Value = -45441

$ws.Range("H53").Value = 1372.3529
$ws.Range("J53").Value = 1930.091
$ws.Range("L53").Value = 1930.091
$ws.Range("N53").Value = -3204.091

$ws.Range("H61").Value = 513.44446
$ws.Range("I61").Value = 540.125
$ws.Range("K61").Value = 1620.375
$ws.Range("M61").Value = -1448.375

$ws.Range("H86").Value = 6646.8335
$ws.Range("I86").Value = 4347.5
$ws.Range("J86").Value = 7796.5
$ws.Range("K86").Value = 4347.5
$ws.Range("L86").Value = 7796.5
$ws.Range("M86").Value = -3224.5
$ws.Range("N86").Value = -10042.5

$ws.Range("H89").Value = 6646.8335
$ws.Range("I89").Value = 4347.5
$ws.Range("J89").Value = 7796.5
$ws.Range("K89").Value = 21737.5
$ws.Range("L89").Value = 38982.5
$ws.Range("M89").Value = -16121.5
$ws.Range("N89").Value = -50214.5

$ws.Range("H100").Value = 4218.5454
$ws.Range("I100").Value = 2650
$ws.Range("J100").Value = 5114.857
$ws.Range("K100").Value = 2650
$ws.Range("L100").Value = 5114.857
$ws.Range("M100").Value = -2109
$ws.Range("N100").Value = -6196.857

$ws.Range("H103").Value = 2440.8572
$ws.Range("I103").Value = 2404
$ws.Range("J103").Value = 2447
$ws.Range("K103").Value = 7212
$ws.Range("L103").Value = 7341
$ws.Range("M103").Value = -6626
$ws.Range("N103").Value = -8513

$ws.Range("H106").Value = 11336.889
$ws.Range("I106").Value = 3394.625
$ws.Range("J106").Value = 17690.7
$ws.Range("K106").Value = 3394.625
$ws.Range("L106").Value = 17690.7
$ws.Range("M106").Value = -2763.625
$ws.Range("N106").Value = -18952.7

$ws.Range("H108").Value = 45952.668
$ws.Range("J108").Value = 45952.668
$ws.Range("L108").Value = 45952.668
$ws.Range("N108").Value = -53632.668

$ws.Range("H109").Value = 52489.332
$ws.Range("J109").Value = 52489.332
$ws.Range("L109").Value = 52489.332
$ws.Range("N109").Value = -55263.332

$ws.Range("H110").Value = 64581.5
$ws.Range("J110").Value = 64581.5
$ws.Range("L110").Value = 64581.5
$ws.Range("N110").Value = -72761.5

$ws.Range("H118").Value = 204.2
$ws.Range("I118").Value = 204.2
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 612.5999999999999
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 1044.4
$ws.Range("N118").ClearContents() | Out-Null

$ws.Range("H132").Value = 1524.7576
$ws.Range("I132").Value = 1350.9286
$ws.Range("K132").Value = 4052.7858
$ws.Range("M132").Value = -1522.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 23530.5
$ws.Range("J42").Value = 23530.5
$ws.Range("L42").Value = 23530.5
$ws.Range("N42").Value = -24502.5

$ws.Range("H45").Value = 71431620
$ws.Range("I45").Value = 90910216
$ws.Range("K45").Value = 90910216
$ws.Range("M45").Value = -90909839

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1269.8
$ws.Range("J80").Value = 971.2857
$ws.Range("L80").Value = 971.2857
$ws.Range("N80").Value = -2967.2857

$ws.Range("H83").Value = 1269.8
$ws.Range("J83").Value = 971.2857
$ws.Range("L83").Value = 4856.4285
$ws.Range("N83").Value = -14840.4285

$ws.Range("H134").Value = 1748.3334
$ws.Range("I134").Value = 1748.3334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5245.0002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2710.0002
$ws.Range("N134").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 7375.8945
$ws.Range("I122").Value = 2269.7
$ws.Range("J122").Value = 13049.444
$ws.Range("K122").Value = 6809.099999999999
$ws.Range("L122").Value = 39148.33199999999
$ws.Range("M122").Value = -4359.099999999999
$ws.Range("N122").Value = -44048.33199999999

$ws.Range("H132").Value = 4279.4
$ws.Range("I132").Value = 4117.933
$ws.Range("J132").Value = 4763.8
$ws.Range("K132").Value = 12353.799
$ws.Range("L132").Value = 14291.4
$ws.Range("M132").Value = -9823.798999999999
$ws.Range("N132").Value = -19351.4

$ws.Range("H134").Value = 2002.9412
$ws.Range("I134").Value = 1265.1538
$ws.Range("J134").Value = 4400.75
$ws.Range("K134").Value = 3795.4614
$ws.Range("L134").Value = 13202.25
$ws.Range("M134").Value = -1260.4614
$ws.Range("N134").Value = -18272.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2199.5
$ws.Range("I81").Value = 2266
$ws.Range("K81").Value = 6798
$ws.Range("M81").Value = -5675

$ws.Range("H84").Value = 2199.5
$ws.Range("I84").Value = 2266
$ws.Range("K84").Value = 20394
$ws.Range("M84").Value = -14778

$ws.Range("H139").Value = 5356.923
$ws.Range("I139").Value = 1626
$ws.Range("K139").Value = 4878
$ws.Range("M139").Value = 262

$ws.Range("H141").Value = 6568.294
$ws.Range("I141").Value = 2557.182
$ws.Range("J141").Value = 13922
$ws.Range("K141").Value = 7671.545999999999
$ws.Range("L141").Value = 41766
$ws.Range("M141").Value = -2491.545999999999
$ws.Range("N141").Value = -52126

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5970.2104
$ws.Range("J80").Value = 6287.769
$ws.Range("L80").Value = 6287.769
$ws.Range("N80").Value = -8283.769

$ws.Range("H83").Value = 5970.2104
$ws.Range("J83").Value = 6287.769
$ws.Range("L83").Value = 31438.845
$ws.Range("N83").Value = -41422.845

$ws.Range("H122").Value = 4675.1904
$ws.Range("I122").Value = 4354.8
$ws.Range("J122").Value = 5476.1665
$ws.Range("K122").Value = 13064.4
$ws.Range("L122").Value = 16428.4995
$ws.Range("M122").Value = -10614.4
$ws.Range("N122").Value = -21328.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2547.353
$ws.Range("I22").Value = 1077.5454
$ws.Range("K22").Value = 1077.5454
$ws.Range("M22").Value = -782.5454

$ws.Range("H27").Value = 2547.353
$ws.Range("I27").Value = 1077.5454
$ws.Range("K27").Value = 1077.5454
$ws.Range("M27").Value = -970.5454

$ws.Range("H40").Value = 7966.3335
$ws.Range("I40").Value = 6777.4287
$ws.Range("K40").Value = 6777.4287
$ws.Range("M40").Value = -6641.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 41126
$ws.Range("J112").Value = 41126
$ws.Range("L112").Value = 41126
$ws.Range("N112").Value = -44080

$ws.Range("H117").Value = 68421
$ws.Range("J117").Value = 68421
$ws.Range("L117").Value = 68421
$ws.Range("N117").Value = -77599

$ws.Range("H126").Value = 3339.182
$ws.Range("J126").Value = 4814.1665
$ws.Range("L126").Value = 14442.4995
$ws.Range("N126").Value = -19382.4995

$ws.Range("H132").Value = 4555.9473
$ws.Range("I132").Value = 4181.1562
$ws.Range("J132").Value = 6554.8335
$ws.Range("K132").Value = 12543.4686
$ws.Range("L132").Value = 19664.5005
$ws.Range("M132").Value = -10013.4686
$ws.Range("N132").Value = -24724.5005
